# Update 2p0. Convention change to support multi-axle vehicles
#
# Renames the single-axle naming convention to an explicit "axle 1" form so
# the workbook's layout can later be extended to multi-axle vehicles:
#   sAxleF      -> sAxle1        (row label, column A, row 5)
#   Body_1Axle  -> Body_Axle1    (class value, column H, row 4)
#
# Applied identically across all three worksheets (Trailer_Elula,
# Trailer_Elula_Unstable, Trailer_Thwala).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A5").Value = "sAxle1"
    $ws.Range("H4").Value = "Body_Axle1"
}

# Restore the active-cell/selection per sheet and re-activate the sheet that
# was originally selected (Trailer_Elula), since selecting cells on the
# other sheets moves the active tab.
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Range("M7").Select() | Out-Null
$ws2.Range("H4").Select() | Out-Null
$ws3.Range("H4").Select() | Out-Null

$ws1.Activate() | Out-Null
